# Commit: Tue, Jul 21, 2020 10:05:35 PM
#
# 1) Swap the table style applied to the cash-flow recap table on slide 16
#    to the new style id.
# 2) Re-colour the presentation's (main) theme from the "Integral" palette
#    to the default "Office Theme" palette -- this is the theme referenced
#    by the slide master / presentation (ppt/theme/theme2.xml), reached via
#    ThemeColorScheme off any slide (they all resolve to the same master
#    colour scheme).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 ------------------------------------------
$s16 = $p.Slides.Item(16)
$tbl = $s16.Shapes.Item(3).Table
$tbl.ApplyStyle("{B8D69958-A70C-4FEF-80D3-33FBAB4A0199}")

# --- 2. Theme colours: Integral -> Office Theme ---------------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
